$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.610.10'
$ws.Range('E2').Value = '  +2.16%  '
$ws.Range('D3').Value = '3.825.88'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '668.60'
$ws.Range('E5').Value = '  +7.62%  '
$ws.Range('D6').Value = '168.91'
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('D7').Value = '3.827.60'
$ws.Range('E7').Value = '  +1.48%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('E10').Value = '  +0.19%  '
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').Value = '7.01'
$ws.Range('E12').Value = '  +5.72%  '
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = '35.78'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').Value = '4.468.42'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '3.814.05'
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('D17').Value = '70.606.60'
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = '17.81'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('E19').Value = '  +1.08%  '
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').Value = '10.77'
$ws.Range('E21').Value = '  +12.30%  '
$ws.Range('D22').Value = '477.15'
$ws.Range('E22').Value = '  +1.89%  '
$ws.Range('D23').Value = '0.714'
$ws.Range('E23').Value = '  +1.90%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '82.98'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '0.0000145'
$ws.Range('E25').Value = '  -2.57%  '
$ws.Range('D26').Value = '12.25'
$ws.Range('E26').Value = '  +2.10%  '
$ws.Range('D27').Value = '10.37'
$ws.Range('E27').Value = '  +3.61%  '
$ws.Range('D28').Value = '2.13'
$ws.Range('E28').Value = '  -1.17%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '3.976.97'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').Value = '2.86'
$ws.Range('E31').Value = '  +7.96%  '
$ws.Range('E32').Value = '  +2.56%  '
$ws.Range('D33').Value = '7.39'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('D34').Value = '29.80'
$ws.Range('E34').Value = '  +3.26%  '
$ws.Range('D35').Value = '0.179'
$ws.Range('E35').Value = '  +13.35%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').Value = '1.01'
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '9.17'
$ws.Range('E37').Value = '  +2.08%  '
$ws.Range('D38').Value = '3.781.08'
$ws.Range('E38').Value = '  +1.59%  '
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('D41').Value = '5.94'
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('D42').Value = '0.974'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '2.08'
$ws.Range('E45').Value = '  +9.44%  '
$ws.Range('D46').Value = '45.42'
$ws.Range('E46').Value = '  +5.80%  '
$ws.Range('D47').Value = '48.78'
$ws.Range('E47').Value = '  +4.30%  '
$ws.Range('D48').Value = '158.54'
$ws.Range('E48').Value = '  +2.81%  '
$ws.Range('B49').Value = 'TheGraph'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D49').Value = '0.301'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '1.43'
$ws.Range('E50').Value = '  +4.53%  '
$ws.Range('D51').Value = '8.50'
$ws.Range('E51').Value = '  +1.21%  '
